# Auto-generated Excel COM-interop script
# Commit message: "Model training severe changes"
# Updates the predicted SHE values (columns A=AVG_SHE, B=L95CI_SHE, C=U95CI_SHE)
# for rows 2-72 on the active worksheet with the newly retrained model's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: @(RowNumber, AVG_SHE, L95CI_SHE, U95CI_SHE)
$data = @(
    @(2, 50.56265333333329, 38.60328000000002, 62.01183999999999),
    @(3, 41.3071052736318, 32.29964000000002, 49.74235999999999),
    @(4, 40.43436, 31.42411999999998, 49.18491999999999),
    @(5, 49.43188059701485, 38.77996000000002, 59.84312),
    @(6, 53.08487044776115, 41.72335999999996, 63.67324000000001),
    @(7, 42.81727900497507, 33.37784, 52.96084000000002),
    @(8, 36.7521643781094, 26.41379999999999, 46.46104000000003),
    @(9, 48.28533134328355, 37.58751999999998, 58.76852000000002),
    @(10, 40.93035940298503, 29.50067999999998, 52.21960000000003),
    @(11, 48.84619940298502, 37.21567999999999, 62.38096000000002),
    @(12, 58.75082208955219, 46.87836000000004, 71.43344000000006),
    @(13, 51.62112696517406, 38.07543999999999, 65.47524000000003),
    @(14, 48.57698149253725, 36.45076, 61.31632),
    @(15, 57.31884955223876, 45.61536000000001, 68.68144000000001),
    @(16, 51.97146646766164, 41.13304, 63.10092),
    @(17, 58.02800298507461, 45.58208, 71.07015999999997),
    @(18, 40.32096656716411, 31.75771999999999, 48.97723999999997),
    @(19, 55.52253174129353, 41.98799999999999, 69.72287999999998),
    @(20, 59.90227263681587, 47.81455999999999, 72.52052),
    @(21, 57.25502467661688, 45.47316000000003, 68.69163999999998),
    @(22, 44.1403824875621, 32.06848000000002, 56.71619999999997),
    @(23, 60.75232716417901, 48.71847999999998, 73.62943999999999),
    @(24, 52.1078051741293, 41.36491999999999, 63.79464),
    @(25, 44.47600199004972, 34.9794, 53.84192000000004),
    @(26, 56.63685233830841, 44.93967999999999, 68.57860000000002),
    @(27, 57.75262089552232, 41.67632, 72.99664000000001),
    @(28, 46.7171460696517, 35.29108000000001, 57.81964000000004),
    @(29, 63.92500577114419, 50.70416, 77.20996000000001),
    @(30, 54.42415223880593, 42.27075999999999, 66.62756000000003),
    @(31, 44.86150587064671, 33.17872, 56.87972000000001),
    @(32, 41.15252796019897, 30.32908, 52.93236000000001),
    @(33, 48.26099761194025, 36.35695999999999, 61.79128000000001),
    @(34, 60.1668302487562, 45.99544, 74.39164000000004),
    @(35, 56.48337850746268, 44.77512000000002, 68.05919999999999),
    @(36, 43.09989791044771, 33.08704, 53.06000000000001),
    @(37, 60.49370427860691, 48.53007999999999, 73.43843999999999),
    @(38, 55.089719800995, 41.77564000000003, 69.65475999999995),
    @(39, 61.42935820895522, 47.86783999999999, 75.19864000000001),
    @(40, 48.48193393034822, 38.29979999999998, 59.14079999999998),
    @(41, 38.82352417910446, 29.73148, 48.24304),
    @(42, 52.47356318407954, 39.39612, 65.67983999999998),
    @(43, 54.35771084577109, 42.66812000000001, 66.6336),
    @(44, 38.44945472636809, 28.65792, 47.83596000000001),
    @(45, 58.87261432835821, 43.45532, 74.34975999999999),
    @(46, 41.28756179104474, 31.11196000000001, 51.36519999999997),
    @(47, 48.47389412935318, 35.05967999999999, 62.09367999999998),
    @(48, 58.76611203980097, 46.63255999999999, 71.64476000000005),
    @(49, 45.65999184079601, 35.83552, 55.16224000000001),
    @(50, 60.19116457711436, 48.20715999999999, 72.81040000000002),
    @(51, 40.83388955223881, 31.69648, 49.64039999999997),
    @(52, 37.82587562189052, 27.91224, 47.25008000000003),
    @(53, 61.07165333333334, 46.67443999999998, 75.40363999999995),
    @(54, 56.89448537313426, 45.20323999999999, 68.55907999999999),
    @(55, 56.77654109452728, 44.93088000000001, 68.71423999999996),
    @(56, 44.92135522388055, 33.14876000000002, 57.14939999999998),
    @(57, 40.99098427860691, 29.72523999999997, 52.02400000000001),
    @(58, 60.62627522388054, 48.42424, 73.19091999999998),
    @(59, 45.12865930348252, 33.34359999999998, 57.38927999999998),
    @(60, 60.36351024875613, 48.39431999999999, 73.14920000000002),
    @(61, 51.74470646766164, 38.35880000000002, 65.46196000000002),
    @(62, 48.90288199004971, 36.95756, 61.5536),
    @(63, 58.32740895522385, 42.94619999999999, 73.46440000000001),
    @(64, 60.55022805970143, 48.81903999999998, 73.36716000000001),
    @(65, 36.90856796019899, 27.07876, 47.74664000000001),
    @(66, 57.15411482587059, 44.62803999999999, 70.73635999999993),
    @(67, 53.4758917412935, 41.26832000000001, 66.5082),
    @(68, 52.88784955223872, 41.89368, 64.55803999999998),
    @(69, 47.26312378109449, 34.49464, 61.40191999999999),
    @(70, 48.92185432835813, 38.66443999999998, 59.66112000000001),
    @(71, 39.76989691542285, 30.62076, 48.87720000000001),
    @(72, 47.43955621890543, 36.52628000000002, 58.07804000000001)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
